$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.252640604972839
$ws.Range("B1").Value = 2.628923177719116
$ws.Range("C1").Value = 4.784129619598389
$ws.Range("D1").Value = 2.249109029769897
$ws.Range("E1").Value = 1.04899263381958
